$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-7 (columns A-T)
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics
$data = @(
    @("FAPs","Rspo3","Lgr5","ECs",3,1,2.507621333333333,7.522864,0.9939780200440224,0.9939780200440224,1,0.3333333333333333,0.08037899999999999,0.241137,0.03101234696192509,0.03101234696192509,0.201560095152,1.814040856368,0.03082559123013256,0.03082559123013256),
    @("FAPs","Rspo3","Lgr5","FAPs",3,1,2.507621333333333,7.522864,0.9939780200440224,0.9939780200440224,3,1,0.8326193333333333,2.497858,0.3212465899369251,0.3212465899369251,2.087894002812444,18.791046025312,0.3193120494113987,0.3193120494113987),
    @("FAPs","Rspo3","Lgr5","sCs",3,1,2.507621333333333,7.522864,0.9939780200440224,0.9939780200440224,2,0.6666666666666666,1.678840333333333,5.036521,0.6477410631011499,0.6477410631011499,4.209895835127111,37.889062516144,0.6438403794024911,0.6438403794024911),
    @("sCs","Rspo3","Lgr5","ECs",1,0.3333333333333333,0.01519233333333333,0.045577,0.006021979955977724,0.006021979955977723,1,0.3333333333333333,0.08037899999999999,0.241137,0.03101234696192509,0.03101234696192509,0.001221144561,0.010990301049,0.0001867557317925396,0.0001867557317925396),
    @("sCs","Rspo3","Lgr5","FAPs",1,0.3333333333333333,0.01519233333333333,0.045577,0.006021979955977724,0.006021979955977723,3,1,0.8326193333333333,2.497858,0.3212465899369251,0.3212465899369251,0.01264943045177778,0.113844874066,0.001934540525526358,0.001934540525526358),
    @("sCs","Rspo3","Lgr5","sCs",1,0.3333333333333333,0.01519233333333333,0.045577,0.006021979955977724,0.006021979955977723,2,0.6666666666666666,1.678840333333333,5.036521,0.6477410631011499,0.6477410631011499,0.02550550195744445,0.229549517617,0.003900683698658827,0.003900683698658826)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex = $colIndex + 1
    }
    $rowIndex = $rowIndex + 1
}
